$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("F3").Value = 2.24
$ws.Range("G3").Value = 3
$ws.Range("H3").Value = 2.92
$ws.Range("I3").Value = 4
$ws.Range("J3").Value = 2.7
$ws.Range("K3").Value = 3.3
$ws.Range("P3").Value = 1.52
$ws.Range("Q3").Value = 2.4

# Row 5
$ws.Range("F5").Value = 1.21
$ws.Range("G5").Value = 1.28
$ws.Range("H5").Value = 11
$ws.Range("I5").Value = 20
$ws.Range("J5").Value = 7
$ws.Range("K5").Value = 9.800000000000001
$ws.Range("N5").Value = 6.6
$ws.Range("O5").Value = 1.13
$ws.Range("P5").Value = 2.94
$ws.Range("Q5").Value = 1.4
$ws.Range("R5").Value = 1.78
$ws.Range("S5").Value = 1.87
$ws.Range("T5").Value = 1.98
$ws.Range("U5").Value = 1.81
$ws.Range("X5").Value = 970
$ws.Range("Y5").Value = 70
$ws.Range("Z5").Value = 190
$ws.Range("AB5").Value = 15
$ws.Range("AC5").Value = 970
$ws.Range("AD5").Value = 65
$ws.Range("AE5").Value = 280
$ws.Range("AF5").Value = 11.5
$ws.Range("AG5").Value = 15
$ws.Range("AH5").Value = 42
$ws.Range("AI5").Value = 190
$ws.Range("AJ5").Value = 12
$ws.Range("AK5").Value = 17
$ws.Range("AL5").Value = 970
$ws.Range("AM5").Value = 190
$ws.Range("AN5").Value = 3.95

# Row 7
$ws.Range("I7").Value = 2.06
$ws.Range("P7").Value = 1.8
$ws.Range("Q7").Value = 2.16
$ws.Range("U7").Value = 1.98

# Row 8
$ws.Range("Q8").Value = 1.51

# Row 9
$ws.Range("O9").Value = 1.26
$ws.Range("P9").Value = 2.12
$ws.Range("Q9").Value = 1.79
$ws.Range("T9").Value = 1.67

# Row 10
$ws.Range("P10").Value = 2.16
$ws.Range("S10").Value = 2.98

# Row 11
$ws.Range("H11").Value = 2.26
$ws.Range("Q11").Value = 1.71
$ws.Range("R11").Value = 1.49

# Row 12
$ws.Range("G12").Value = 1.61
$ws.Range("P12").Value = 2.2
$ws.Range("T12").Value = 1.86
$ws.Range("U12").Value = 2.06
$ws.Range("X12").Value = 24

# Row 13
$ws.Range("G13").Value = 2.72
$ws.Range("U13").Value = 2.38
$ws.Range("X13").Value = 23
$ws.Range("AA13").Value = 48
$ws.Range("AF13").Value = 19.5
$ws.Range("AI13").Value = 42
$ws.Range("AJ13").Value = 46
$ws.Range("AM13").Value = 85

# Row 14
$ws.Range("H14").Value = 1.8
$ws.Range("S14").Value = 2.64
$ws.Range("U14").Value = 2.34
$ws.Range("AO14").Value = 9

# Row 15
$ws.Range("J15").Value = 5.3
$ws.Range("Q15").Value = 1.47
$ws.Range("R15").Value = 1.74
$ws.Range("AN15").Value = 5

# Row 16
$ws.Range("F16").Value = 4
$ws.Range("G16").Value = 5.2
$ws.Range("H16").Value = 2.28
$ws.Range("I16").Value = 2.78
$ws.Range("J16").Value = 2.48
$ws.Range("K16").Value = 3.05
$ws.Range("P16").Value = 1.39

# Row 17
$ws.Range("O17").Value = 1.31
$ws.Range("U17").Value = 1.96

# Row 18
$ws.Range("N18").Value = 4
$ws.Range("O18").Value = 1.27
$ws.Range("P18").Value = 2.04
$ws.Range("Q18").Value = 1.8
$ws.Range("R18").Value = 1.42
$ws.Range("S18").Value = 3
$ws.Range("T18").Value = 1.72
$ws.Range("U18").Value = 2.16
$ws.Range("X18").Value = 17.5
$ws.Range("Y18").Value = 17.5
$ws.Range("Z18").Value = 34
$ws.Range("AA18").Value = 110
$ws.Range("AB18").Value = 10.5
$ws.Range("AC18").Value = 9.199999999999999
$ws.Range("AD18").Value = 18
$ws.Range("AF18").Value = 13
$ws.Range("AG18").Value = 11
$ws.Range("AH18").Value = 18
$ws.Range("AI18").Value = 60
$ws.Range("AJ18").Value = 23
$ws.Range("AK18").Value = 20
$ws.Range("AL18").Value = 34
$ws.Range("AM18").Value = 110
$ws.Range("AN18").Value = 12

# Row 19
$ws.Range("P19").Value = 2
$ws.Range("S19").Value = 3.35
$ws.Range("T19").Value = 1.74
$ws.Range("AM19").Value = 85

# Row 20
$ws.Range("F20").Value = 2.8
$ws.Range("H20").Value = 2.6
$ws.Range("I20").Value = 2.68
$ws.Range("J20").Value = 3.6
$ws.Range("O20").Value = 1.29
$ws.Range("Q20").Value = 1.86
$ws.Range("S20").Value = 3.15
$ws.Range("T20").Value = 1.7
$ws.Range("U20").Value = 2.26
$ws.Range("AC20").Value = 8.4

# Row 21
$ws.Range("P21").Value = 1.95

# Row 22
$ws.Range("G22").Value = 1.98
$ws.Range("H22").Value = 4
$ws.Range("I22").Value = 4.5
$ws.Range("K22").Value = 4.4
$ws.Range("P22").Value = 2.3

# Row 23
$ws.Range("J23").Value = 3.8
$ws.Range("K23").Value = 3.9

# Row 24
$ws.Range("J24").Value = 3.45
$ws.Range("K24").Value = 3.55
$ws.Range("N24").Value = 3.7
$ws.Range("Q24").Value = 1.95
$ws.Range("R24").Value = 1.36

# Row 25
$ws.Range("I25").Value = 2.92
$ws.Range("S25").Value = 2.48
